$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.370.73"
$ws.Range("E2").Value = "  +9.57%  "

$ws.Range("D3").Value = "1.783.77"
$ws.Range("E3").Value = "  +6.50%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'337.67"
$ws.Range("E5").Value = "  +2.36%  "

$ws.Range("D6").Value = "'0.9964"
$ws.Range("E6").Value = "  -0.31%  "

$ws.Range("D7").Value = "'0.3796"
$ws.Range("E7").Value = "  +3.83%  "

$ws.Range("B8").Value = "OKB"
$ws.Range("C8").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D8").Value = "'49.31"
$ws.Range("E8").Value = "  +4.70%  "

$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "'0.3468"
$ws.Range("E9").Value = "  +7.09%  "

$ws.Range("D10").Value = "'1.215"
$ws.Range("E10").Value = "  +6.03%  "

$ws.Range("D11").Value = "'0.07653"
$ws.Range("E11").Value = "  +6.75%  "

$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  +0.06%  "

$ws.Range("D13").Value = "'6.562"
$ws.Range("E13").Value = "  +7.85%  "

$ws.Range("D14").Value = "'21.34"
$ws.Range("E14").Value = "  +8.47%  "

$ws.Range("D15").Value = "'7.212"
$ws.Range("E15").Value = "  +8.23%  "

$ws.Range("D16").Value = "1.783.91"
$ws.Range("E16").Value = "  +6.72%  "

$ws.Range("D17").Value = "'0.00001113"
$ws.Range("E17").Value = "  +6.24%  "

$ws.Range("D18").Value = "'0.06759"
$ws.Range("E18").Value = "  +3.47%  "

$ws.Range("D19").Value = "'84.80"
$ws.Range("E19").Value = "  +7.50%  "

$ws.Range("D20").Value = "'0.9971"
$ws.Range("E20").Value = "  -0.21%  "

$ws.Range("D21").Value = "'17.56"
$ws.Range("E21").Value = "  +10.91%  "

$ws.Range("D22").Value = "'6.372"
$ws.Range("E22").Value = "  +7.85%  "

$ws.Range("D23").Value = "'13.10"
$ws.Range("E23").Value = "  +1.32%  "

$ws.Range("D24").Value = "27.404.32"
$ws.Range("E24").Value = "  +9.66%  "

$ws.Range("D25").Value = "'2.471"
$ws.Range("E25").Value = "  +1.41%  "

$ws.Range("D26").Value = "'1.510"
$ws.Range("E26").Value = "  +26.12%  "

$ws.Range("D27").Value = "'2.522"
$ws.Range("E27").Value = "  +6.26%  "

$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "'153.13"
$ws.Range("E28").Value = "  +2.56%  "

$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'20.14"
$ws.Range("E29").Value = "  +7.58%  "

$ws.Range("D30").Value = "1.979.51"
$ws.Range("E30").Value = "  +6.54%  "

$ws.Range("D31").Value = "'135.20"
$ws.Range("E31").Value = "  +7.33%  "

$ws.Range("D32").Value = "'4.104"
$ws.Range("E32").Value = "  +0.49%  "

$ws.Range("D33").Value = "'6.264"
$ws.Range("E33").Value = "  +8.16%  "

$ws.Range("D34").Value = "'0.08724"
$ws.Range("E34").Value = "  +3.08%  "

$ws.Range("D35").Value = "'13.41"
$ws.Range("E35").Value = "  +8.86%  "

$ws.Range("D36").Value = "'1.718"
$ws.Range("E36").Value = "  +2.83%  "

$ws.Range("D37").Value = "'5.588"
$ws.Range("E37").Value = "  +8.50%  "

$ws.Range("D38").Value = "'0.02410"
$ws.Range("E38").Value = "  +8.26%  "

$ws.Range("D39").Value = "'0.2242"
$ws.Range("E39").Value = "  +7.32%  "

$ws.Range("E40").Value = "  +6.26%  "

$ws.Range("D41").Value = "'0.6760"
$ws.Range("E41").Value = "  +13.43%  "

$ws.Range("D42").Value = "'8.747"
$ws.Range("E42").Value = "  +5.98%  "

$ws.Range("D43").Value = "'1.234"
$ws.Range("E43").Value = "  -0.07%  "

$ws.Range("D44").Value = "'14.73"
$ws.Range("E44").Value = "  +7.79%  "

$ws.Range("D45").Value = "'0.6425"
$ws.Range("E45").Value = "  +12.21%  "

$ws.Range("D46").Value = "'0.9961"
$ws.Range("E46").Value = "  -0.27%  "

$ws.Range("D47").Value = "'3.964"
$ws.Range("E47").Value = "  +3.57%  "

$ws.Range("D48").Value = "'2.141"
$ws.Range("E48").Value = "  +9.09%  "

$ws.Range("D49").Value = "'130.79"
$ws.Range("E49").Value = "  +5.46%  "

$ws.Range("D50").Value = "'0.07343"
$ws.Range("E50").Value = "  +4.73%  "

$ws.Range("D51").Value = "'79.92"
$ws.Range("E51").Value = "  +7.18%  "
